$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update measurement data for rows 17-23
$ws.Range("E17").Value = 10015
$ws.Range("F17").Value = 19.181827259111301
$ws.Range("G17").Value = 0.089465801298052899

$ws.Range("E18").Value = 10000
$ws.Range("F18").Value = 18.872800000000002
$ws.Range("G18").Value = 0.091499999999999998

$ws.Range("E19").Value = 10000
$ws.Range("F19").Value = 17.988299999999999
$ws.Range("G19").Value = 0.087300000000000003

$ws.Range("E20").Value = 10000
$ws.Range("F20").Value = 17.8934
$ws.Range("G20").Value = 0.085699999999999998

$ws.Range("E21").Value = 10000
$ws.Range("F21").Value = 17.716999999999999
$ws.Range("G21").Value = 0.086300000000000002

$ws.Range("E22").Value = 10000
$ws.Range("F22").Value = 16.899100000000001
$ws.Range("G22").Value = 0.0872

$ws.Range("E23").Value = 10062
$ws.Range("F23").Value = 4.0263367123832197
$ws.Range("G23").Value = 0.090141125024845906

# Update measurement data for rows 27-33
$ws.Range("E27").Value = 10015
$ws.Range("F27").Value = 16.4982526210683
$ws.Range("G27").Value = 1.03764353469795

$ws.Range("E28").Value = 10000
$ws.Range("F28").Value = 17.379899999999999
$ws.Range("G28").Value = 1.5185999999999999

$ws.Range("E29").Value = 10000
$ws.Range("F29").Value = 17.5656
$ws.Range("G29").Value = 1.8375999999999999

$ws.Range("E30").Value = 10391
$ws.Range("F30").Value = 15.6894427870272
$ws.Range("G30").Value = 1.6803964969685301

$ws.Range("E31").Value = 10000
$ws.Range("F31").Value = 16.2425
$ws.Range("G31").Value = 1.6834

$ws.Range("E32").Value = 10000
$ws.Range("F32").Value = 15.190300000000001
$ws.Range("G32").Value = 1.9419999999999999

$ws.Range("E33").Value = 10063
$ws.Range("F33").Value = 3.5943555599721702
$ws.Range("G33").Value = 1.8674351585014399

# Update the selected cell to G33 as recorded in the sheet view
$ws.Range("G33").Select()
